$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the empty spacer paragraphs (style BodyText / FirstParagraph,
#    no text) that sit just before each w:bookmarkStart/w:bookmarkEnd
#    section marker throughout the document - except the very last one,
#    which sits just before the "NCSC" section's bookmarkEnd and instead
#    gets turned into a horizontal-rule paragraph (step 2).
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
$emptyIdx = @()
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.Trim() -eq "") {
        $emptyIdx += $i
    }
}

# The last empty paragraph in document order is the one right before the
# "NCSC" section's closing bookmark - keep it (it becomes the horizontal
# rule below). Delete every other empty paragraph, walking backwards so
# earlier indices stay valid while we delete.
$lastEmptyOriginal = $emptyIdx[$emptyIdx.Count - 1]
$numDeletedBeforeLast = $emptyIdx.Count - 1
for ($j = $emptyIdx.Count - 2; $j -ge 0; $j--) {
    $d.Paragraphs($emptyIdx[$j]).Range.Delete()
}

# ------------------------------------------------------------------
# 2. Turn the remaining empty paragraph (just before the NCSC section's
#    bookmarkEnd) into a horizontal-rule paragraph: <w:r><w:pict><v:rect
#    style="width:0;height:1.5pt" o:hralign="center" o:hrstd="t" o:hr="t"/>
#    </w:pict></w:r>
# ------------------------------------------------------------------
$hrIdx = $lastEmptyOriginal - $numDeletedBeforeLast
$hrParagraph = $d.Paragraphs($hrIdx)
$hrXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:o="urn:schemas-microsoft-com:office:office">' +
  '<w:body><w:p><w:r><w:pict>' +
  '<rect xmlns="urn:schemas-microsoft-com:vml" style="width:0;height:1.5pt" o:hralign="center" o:hrstd="t" o:hr="t"/>' +
  '</w:pict></w:r></w:p></w:body></w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'
$hrParagraph.Range.InsertXML($hrXml)

# ------------------------------------------------------------------
# 3. The closing "Feedback" paragraph switches from the BlockText style
#    to FirstParagraph.
# ------------------------------------------------------------------
$feedbackParagraph = $d.Paragraphs($d.Paragraphs.Count)
$feedbackParagraph.Style = "First Paragraph"
